$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 7926
$ws.Range("F11").Value = 430
$ws.Range("F12").Value = 1705
$ws.Range("F13").Value = 72
$ws.Range("F14").Value = 1108
$ws.Range("F18").Value = 8570
$ws.Range("F19").Value = 209
$ws.Range("F33").Value = 1056
$ws.Range("F37").Value = 3577
$ws.Range("F46").Value = 55
$ws.Range("F47").Value = 114
$ws.Range("F48").Value = 30
$ws.Range("F49").Value = 2419

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 2251
$ws.Range("F9").Value = 9078

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7926
$ws.Range("F6").Value = 2251
$ws.Range("F12").Value = 430
$ws.Range("F13").Value = 1705
$ws.Range("F14").Value = 72
$ws.Range("F15").Value = 1108
$ws.Range("F17").Value = 8570
$ws.Range("F18").Value = 209
$ws.Range("F30").Value = 1056
$ws.Range("F36").Value = 3577
$ws.Range("F44").Value = 55
$ws.Range("F46").Value = 30
$ws.Range("F48").Value = 2419
